$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")
$ws.Activate()

# Remove the two trailing blank rows (27 and 28) entirely so the new
# content lands cleanly on rows 28 and 29 without leftover row metadata.
$ws.Range("A27:A28").EntireRow.Delete()

# New "Transação de Serviços" section title (row 28), styled the same way
# as the other section titles (row 1 / row 23), then merged across A:B.
$ws.Range("A1:B1").Copy()
$ws.Range("A28:B28").PasteSpecial(-4122)
$ws.Range("A28").Value = "Transação de Serviços"
$ws.Range("A28:B28").Merge()

# New data row (row 29)
$ws.Range("A29").Value = 20
$ws.Range("A29").HorizontalAlignment = -4108
$ws.Range("B29").Value = "Produtos inválidos"
$ws.Range("C29").Value = "Códigos de Produtos informados não encontrados na base ou Produtos não vinculados ao Estabelecimentos"

# Restore the default row height (12.8) on rows whose height had drifted
# to 12.1 in the source file.
$ws.Rows.Item(1).AutoFit()
$ws.Rows.Item(20).AutoFit()
$ws.Rows.Item(22).AutoFit()
$ws.Rows.Item(23).AutoFit()
$ws.Rows.Item(24).AutoFit()

$ws.Application.Goto($ws.Range("A7"), $false)
$ws.Range("C30").Select()
